$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Master")

$nisn = @(
    "0054497273",
    "0054497274",
    "0054497275",
    "0054497276",
    "0054497277",
    "0054497278",
    "0054497279",
    "0054497280",
    "0054497281",
    "0054497282"
)

$startRow = 4
$startDate = 44544

for ($i = 0; $i -lt $nisn.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $nisn[$i]
    $ws.Range("B$row").Value = $row
    $ws.Range("C$row").Value = $startDate + $i
    $ws.Range("D$row").Value = "Akutuh"
}

$ws.Range("A5:D13").Select()
